$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.652.15"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.44"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.43"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4662"
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07137"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9040"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07747"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.46"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.91"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.267"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.347"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.67"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008556"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.691.94"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.015"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.915"
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.72"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.94"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.980"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.01"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.847"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08807"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.150"
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.824"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.167"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7409"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.449"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.078"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01926"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.929"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05138"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.888"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5076"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.042"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4675"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.989"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.33"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.570"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06062"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.16"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.89"
